# Update 想去人数 (F) and 最低票价 (G) figures on the "展览" and "全部类型"
# sheets to match the freshly generated data snapshot.
#
# Row map (same on both sheets):
#   row 2  : F 730 -> 729 ; G 68 (number) -> "不可售" (inline string)
#   row 3  : F 592 -> 594
#   row 4  : F 553 -> 556
#   row 7  : F 74  -> 76
#   row 11 : F 4780 -> 4788
#   row 12 : F 4527 -> 4531
#   row 16 : F 31  -> 32
#   row 17 : F 164 -> 166

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 729
    $ws.Range("G2").Value = "不可售"

    $ws.Range("F3").Value = 594

    $ws.Range("F4").Value = 556

    $ws.Range("F7").Value = 76

    $ws.Range("F11").Value = 4788

    $ws.Range("F12").Value = 4531

    $ws.Range("F16").Value = 32

    $ws.Range("F17").Value = 166
}
